# Updated symbol list on Wed Dec 21 03:50:48 UTC 2022 with GitHub Actions
# Applies the price/coin-ranking refresh described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) refreshes that stay on the same row ---
$ws.Range("D2").Value  = "'248.83"
$ws.Range("D3").Value  = "'22.57"
$ws.Range("D4").Value  = "'5.395"
$ws.Range("D5").Value  = "'0.05715"
$ws.Range("D6").Value  = "'3.413"
$ws.Range("D7").Value  = "'6.331"
$ws.Range("D9").Value  = "'0.9238"

# --- Rows 10-19: coin ranking list shifted up by one (WazirX dropped to the
#     bottom/replaced), each row taking on the coin/link/rank-label of the
#     row below it, along with refreshed prices ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01135"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1422"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07448"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03117"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03029"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09354"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.725"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001589"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04754"
$ws.Range("E18").Value = "17CoinExTokenCET"

$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").Value = "'0.01829"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"

# --- Further price (column D) refreshes further down the table ---
$ws.Range("D20").Value = "'0.006433"
$ws.Range("D21").Value = "'0.005010"
$ws.Range("D22").Value = "'0.001026"
$ws.Range("D24").Value = "'3.699"
$ws.Range("D25").Value = "'2.164"
$ws.Range("D26").Value = "'0.3300"
$ws.Range("D27").Value = "'0.1307"
$ws.Range("D40").Value = "'0.03988"
$ws.Range("D41").Value = "'0.006895"
$ws.Range("D43").Value = "'0.002712"

# --- Label tweak: CoinbaseStockToken now also flagged as the day's worst performer ---
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
